$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1152.5352  # H15: 12988023 -> 1152.5352
$ws.Cells.Item(15, 9).Value = 1152.5352  # I15: 12988023 -> 1152.5352
$ws.Cells.Item(15, 11).Value = 3457.6056  # K15: 38964069 -> 3457.6056
$ws.Cells.Item(15, 13).Value = -3288.6056  # M15: -38963900 -> -3288.6056
$ws.Cells.Item(33, 8).Value = 1062.9714  # H33: 1065.2572 -> 1062.9714
$ws.Cells.Item(33, 9).Value = 307.13794  # I33: 302.9 -> 307.13794
$ws.Cells.Item(33, 10).Value = 4716.1665  # J33: 5639.4 -> 4716.1665
$ws.Cells.Item(33, 11).Value = 307.13794  # K33: 302.9 -> 307.13794
$ws.Cells.Item(33, 12).Value = 4716.1665  # L33: 5639.4 -> 4716.1665
$ws.Cells.Item(33, 13).Value = -78.13794000000001  # M33: -73.89999999999998 -> -78.13794000000001
$ws.Cells.Item(33, 14).Value = -5174.1665  # N33: -6097.4 -> -5174.1665
$ws.Cells.Item(129, 8).Value = 5717.16  # H129: 2473.7542 -> 5717.16
$ws.Cells.Item(129, 10).Value = 2999.0833  # J129: 915.8125 -> 2999.0833
$ws.Cells.Item(129, 12).Value = 8997.249899999999  # L129: 2747.4375 -> 8997.249899999999
$ws.Cells.Item(129, 14).Value = -18997.2499  # N129: -12747.4375 -> -18997.2499
$ws.Cells.Item(132, 8).Value = 6585294.5  # H132: 7582951.5 -> 6585294.5
$ws.Cells.Item(132, 9).Value = 6950808  # I132: 7819571.5 -> 6950808
$ws.Cells.Item(132, 10).Value = 6058.5  # J132: 11111 -> 6058.5
$ws.Cells.Item(132, 11).Value = 20852424  # K132: 23458714.5 -> 20852424
$ws.Cells.Item(132, 12).Value = 18175.5  # L132: 33333 -> 18175.5
$ws.Cells.Item(132, 13).Value = -20849894  # M132: -23456184.5 -> -20849894
$ws.Cells.Item(132, 14).Value = -23235.5  # N132: -38393 -> -23235.5
$ws.Cells.Item(135, 8).Value = 510.88135  # H135: 523.0714 -> 510.88135
$ws.Cells.Item(135, 9).Value = 428.48148  # I135: 442.07693 -> 428.48148
$ws.Cells.Item(135, 10).Value = 1400.8  # J135: 1576 -> 1400.8
$ws.Cells.Item(135, 11).Value = 3856.33332  # K135: 3978.69237 -> 3856.33332
$ws.Cells.Item(135, 12).Value = 12607.2  # L135: 14184 -> 12607.2
$ws.Cells.Item(135, 13).Value = -1321.33332  # M135: -1443.69237 -> -1321.33332
$ws.Cells.Item(135, 14).Value = -17677.2  # N135: -19254 -> -17677.2
$ws.Cells.Item(137, 8).Value = 1432.5135  # H137: 1305.4103 -> 1432.5135
$ws.Cells.Item(137, 9).Value = 1182.4706  # I137: 1135.4324 -> 1182.4706
$ws.Cells.Item(137, 10).Value = 4266.3335  # J137: 4450 -> 4266.3335
$ws.Cells.Item(137, 11).Value = 3547.4118  # K137: 3406.2972 -> 3547.4118
$ws.Cells.Item(137, 12).Value = 12799.0005  # L137: 13350 -> 12799.0005
$ws.Cells.Item(137, 13).Value = -997.4118000000003  # M137: -856.2972 -> -997.4118000000003
$ws.Cells.Item(137, 14).Value = -17899.0005  # N137: -18450 -> -17899.0005
$ws.Cells.Item(141, 8).Value = 1918.6444  # H141: 1868.3489 -> 1918.6444
$ws.Cells.Item(141, 9).Value = 1771.683  # I141: 1771.8049 -> 1771.683
$ws.Cells.Item(141, 10).Value = 3425  # J141: 3847.5 -> 3425
$ws.Cells.Item(141, 11).Value = 5315.049  # K141: 5315.4147 -> 5315.049
$ws.Cells.Item(141, 12).Value = 10275  # L141: 11542.5 -> 10275
$ws.Cells.Item(141, 13).Value = -135.049  # M141: -135.4147000000003 -> -135.049
$ws.Cells.Item(141, 14).Value = -20635  # N141: -21902.5 -> -20635

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19042.566  # H32: 20251.988 -> 19042.566
$ws.Cells.Item(32, 9).Value = 4184.5454  # I32: 4469.0366 -> 4184.5454
$ws.Cells.Item(32, 11).Value = 4184.5454  # K32: 4469.0366 -> 4184.5454
$ws.Cells.Item(32, 13).Value = -3897.5454  # M32: -4182.0366 -> -3897.5454
$ws.Cells.Item(61, 8).Value = 2246.9512  # H61: 1973.36 -> 2246.9512
$ws.Cells.Item(61, 9).Value = 1376.5625  # I61: 1171.92 -> 1376.5625
$ws.Cells.Item(61, 10).Value = 2804  # J61: 2774.8 -> 2804
$ws.Cells.Item(61, 11).Value = 1376.5625  # K61: 1171.92 -> 1376.5625
$ws.Cells.Item(61, 12).Value = 2804  # L61: 2774.8 -> 2804
$ws.Cells.Item(61, 13).Value = -1164.5625  # M61: -959.9200000000001 -> -1164.5625
$ws.Cells.Item(61, 14).Value = -3228  # N61: -3198.8 -> -3228
$ws.Cells.Item(74, 8).Value = 800.4516  # H74: 793.94116 -> 800.4516
$ws.Cells.Item(74, 9).Value = 798.61536  # I74: 816.5599999999999 -> 798.61536
$ws.Cells.Item(74, 10).Value = 810  # J74: 731.1111 -> 810
$ws.Cells.Item(74, 11).Value = 798.61536  # K74: 816.5599999999999 -> 798.61536
$ws.Cells.Item(74, 12).Value = 810  # L74: 731.1111 -> 810
$ws.Cells.Item(74, 13).Value = 75.38463999999999  # M74: 57.44000000000005 -> 75.38463999999999
$ws.Cells.Item(74, 14).Value = -2558  # N74: -2479.1111 -> -2558
$ws.Cells.Item(77, 8).Value = 800.4516  # H77: 793.94116 -> 800.4516
$ws.Cells.Item(77, 9).Value = 798.61536  # I77: 816.5599999999999 -> 798.61536
$ws.Cells.Item(77, 10).Value = 810  # J77: 731.1111 -> 810
$ws.Cells.Item(77, 11).Value = 3993.0768  # K77: 4082.8 -> 3993.0768
$ws.Cells.Item(77, 12).Value = 4050  # L77: 3655.5555 -> 4050
$ws.Cells.Item(77, 13).Value = 374.9232000000002  # M77: 285.2000000000003 -> 374.9232000000002
$ws.Cells.Item(77, 14).Value = -12786  # N77: -12391.5555 -> -12786
$ws.Cells.Item(132, 8).Value = 2426.1162  # H132: 4075 -> 2426.1162
$ws.Cells.Item(132, 9).Value = 2088.8333  # I132: 4546.778 -> 2088.8333
$ws.Cells.Item(132, 10).Value = 3204.4614  # J132: 3367.3333 -> 3204.4614
$ws.Cells.Item(132, 11).Value = 6266.499899999999  # K132: 13640.334 -> 6266.499899999999
$ws.Cells.Item(132, 12).Value = 9613.3842  # L132: 10101.9999 -> 9613.3842
$ws.Cells.Item(132, 13).Value = -3736.499899999999  # M132: -11110.334 -> -3736.499899999999
$ws.Cells.Item(132, 14).Value = -14673.3842  # N132: -15161.9999 -> -14673.3842
$ws.Cells.Item(136, 8).Value = 2246.9512  # H136: 1973.36 -> 2246.9512
$ws.Cells.Item(136, 9).Value = 1376.5625  # I136: 1171.92 -> 1376.5625
$ws.Cells.Item(136, 10).Value = 2804  # J136: 2774.8 -> 2804
$ws.Cells.Item(136, 11).Value = 4129.6875  # K136: 3515.76 -> 4129.6875
$ws.Cells.Item(136, 12).Value = 8412  # L136: 8324.400000000001 -> 8412
$ws.Cells.Item(136, 13).Value = -1579.6875  # M136: -965.7600000000002 -> -1579.6875
$ws.Cells.Item(136, 14).Value = -13512  # N136: -13424.4 -> -13512

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 43782.777  # H86: 42375.035 -> 43782.777
$ws.Cells.Item(86, 9).Value = 58060.25  # I86: 53010.09 -> 58060.25
$ws.Cells.Item(86, 10).Value = 2990  # J86: 3379.8333 -> 2990
$ws.Cells.Item(86, 11).Value = 58060.25  # K86: 53010.09 -> 58060.25
$ws.Cells.Item(86, 12).Value = 2990  # L86: 3379.8333 -> 2990
$ws.Cells.Item(86, 13).Value = -56937.25  # M86: -51887.09 -> -56937.25
$ws.Cells.Item(86, 14).Value = -5236  # N86: -5625.8333 -> -5236
$ws.Cells.Item(89, 8).Value = 43782.777  # H89: 42375.035 -> 43782.777
$ws.Cells.Item(89, 9).Value = 58060.25  # I89: 53010.09 -> 58060.25
$ws.Cells.Item(89, 10).Value = 2990  # J89: 3379.8333 -> 2990
$ws.Cells.Item(89, 11).Value = 290301.25  # K89: 265050.45 -> 290301.25
$ws.Cells.Item(89, 12).Value = 14950  # L89: 16899.1665 -> 14950
$ws.Cells.Item(89, 13).Value = -284685.25  # M89: -259434.45 -> -284685.25
$ws.Cells.Item(89, 14).Value = -26182  # N89: -28131.1665 -> -26182
$ws.Cells.Item(134, 8).Value = 7370.263  # H134: 6910.7085 -> 7370.263
$ws.Cells.Item(134, 9).Value = 7202.3335  # I134: 6442.9 -> 7202.3335
$ws.Cells.Item(134, 10).Value = 8000  # J134: 9249.75 -> 8000
$ws.Cells.Item(134, 11).Value = 21607.0005  # K134: 19328.7 -> 21607.0005
$ws.Cells.Item(134, 12).Value = 24000  # L134: 27749.25 -> 24000
$ws.Cells.Item(134, 13).Value = -19072.0005  # M134: -16793.7 -> -19072.0005
$ws.Cells.Item(134, 14).Value = -29070  # N134: -32819.25 -> -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1748.75  # H16: 1277.7646 -> 1748.75
$ws.Cells.Item(16, 9).Value = 1250  # I16: 1087.6428 -> 1250
$ws.Cells.Item(16, 10).Value = 2247.5  # J16: 2165 -> 2247.5
$ws.Cells.Item(16, 11).Value = 1250  # K16: 1087.6428 -> 1250
$ws.Cells.Item(16, 12).Value = 2247.5  # L16: 2165 -> 2247.5
$ws.Cells.Item(16, 13).Value = -963  # M16: -800.6428000000001 -> -963
$ws.Cells.Item(16, 14).Value = -2821.5  # N16: -2739 -> -2821.5
$ws.Cells.Item(58, 8).Value = 1044.3773  # H58: 568 -> 1044.3773
$ws.Cells.Item(58, 9).Value = 924.16  # I58: 522 -> 924.16
$ws.Cells.Item(58, 10).Value = 3048  # J58: 614 -> 3048
$ws.Cells.Item(58, 11).Value = 924.16  # K58: 522 -> 924.16
$ws.Cells.Item(58, 12).Value = 3048  # L58: 614 -> 3048
$ws.Cells.Item(58, 13).Value = -721.16  # M58: -319 -> -721.16
$ws.Cells.Item(58, 14).Value = -3454  # N58: -1020 -> -3454
$ws.Cells.Item(99, 8).Value = 12530.588  # H99: 12472.471 -> 12530.588
$ws.Cells.Item(99, 9).Value = 5774.2  # I99: 5576.6 -> 5774.2
$ws.Cells.Item(99, 11).Value = 5774.2  # K99: 5576.6 -> 5774.2
$ws.Cells.Item(99, 13).Value = -4276.2  # M99: -4078.6 -> -4276.2
$ws.Cells.Item(113, 8).Value = 1748.75  # H113: 1277.7646 -> 1748.75
$ws.Cells.Item(113, 9).Value = 1250  # I113: 1087.6428 -> 1250
$ws.Cells.Item(113, 10).Value = 2247.5  # J113: 2165 -> 2247.5
$ws.Cells.Item(113, 11).Value = 1250  # K113: 1087.6428 -> 1250
$ws.Cells.Item(113, 12).Value = 2247.5  # L113: 2165 -> 2247.5
$ws.Cells.Item(113, 13).Value = 920  # M113: 1082.3572 -> 920
$ws.Cells.Item(113, 14).Value = -6587.5  # N113: -6505 -> -6587.5
$ws.Cells.Item(126, 8).Value = 12530.588  # H126: 12472.471 -> 12530.588
$ws.Cells.Item(126, 9).Value = 5774.2  # I126: 5576.6 -> 5774.2
$ws.Cells.Item(126, 11).Value = 17322.6  # K126: 16729.8 -> 17322.6
$ws.Cells.Item(126, 13).Value = -14852.6  # M126: -14259.8 -> -14852.6
$ws.Cells.Item(134, 8).Value = 1274.1936  # H134: 1289.5 -> 1274.1936
$ws.Cells.Item(134, 9).Value = 1173.421  # I134: 1168.1052 -> 1173.421
$ws.Cells.Item(134, 10).Value = 1433.75  # J134: 1499.1818 -> 1433.75
$ws.Cells.Item(134, 11).Value = 3520.263  # K134: 3504.3156 -> 3520.263
$ws.Cells.Item(134, 12).Value = 4301.25  # L134: 4497.5454 -> 4301.25
$ws.Cells.Item(134, 13).Value = -985.2629999999999  # M134: -969.3155999999999 -> -985.2629999999999
$ws.Cells.Item(134, 14).Value = -9371.25  # N134: -9567.545399999999 -> -9371.25
$ws.Cells.Item(136, 8).Value = 1044.3773  # H136: 568 -> 1044.3773
$ws.Cells.Item(136, 9).Value = 924.16  # I136: 522 -> 924.16
$ws.Cells.Item(136, 10).Value = 3048  # J136: 614 -> 3048
$ws.Cells.Item(136, 11).Value = 2772.48  # K136: 1566 -> 2772.48
$ws.Cells.Item(136, 12).Value = 9144  # L136: 1842 -> 9144
$ws.Cells.Item(136, 13).Value = -222.48  # M136: 984 -> -222.48
$ws.Cells.Item(136, 14).Value = -14244  # N136: -6942 -> -14244

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 610740.2  # H37: 578143.9 -> 610740.2
$ws.Cells.Item(37, 10).Value = 610740.2  # J37: 578143.9 -> 610740.2
$ws.Cells.Item(37, 12).Value = 1832220.6  # L37: 1734431.7 -> 1832220.6
$ws.Cells.Item(37, 14).Value = -1832444.6  # N37: -1734655.7 -> -1832444.6
$ws.Cells.Item(131, 8).Value = 8692.064  # H131: 8689.598 -> 8692.064
$ws.Cells.Item(131, 10).Value = 8849.187  # J131: 8846.653 -> 8849.187
$ws.Cells.Item(131, 12).Value = 26547.561  # L131: 26539.959 -> 26547.561
$ws.Cells.Item(131, 14).Value = -36627.561  # N131: -36619.959 -> -36627.561
$ws.Cells.Item(139, 8).Value = 2096.639  # H139: 1818.421 -> 2096.639
$ws.Cells.Item(139, 9).Value = 1208.2778  # I139: 1018.5238 -> 1208.2778
$ws.Cells.Item(139, 10).Value = 2985  # J139: 2806.5293 -> 2985
$ws.Cells.Item(139, 11).Value = 3624.8334  # K139: 3055.5714 -> 3624.8334
$ws.Cells.Item(139, 12).Value = 8955  # L139: 8419.5879 -> 8955
$ws.Cells.Item(139, 13).Value = 1515.1666  # M139: 2084.4286 -> 1515.1666
$ws.Cells.Item(139, 14).Value = -19235  # N139: -18699.5879 -> -19235

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 4288085.5  # H7: 4583433.5 -> 4288085.5
$ws.Cells.Item(7, 9).Value = 4288085.5  # I7: 4583433.5 -> 4288085.5
$ws.Cells.Item(7, 11).Value = 4288085.5  # K7: 4583433.5 -> 4288085.5
$ws.Cells.Item(7, 13).Value = -4287973.5  # M7: -4583321.5 -> -4287973.5
$ws.Cells.Item(8, 8).Value = 4288085.5  # H8: 4583433.5 -> 4288085.5
$ws.Cells.Item(8, 9).Value = 4288085.5  # I8: 4583433.5 -> 4288085.5
$ws.Cells.Item(8, 11).Value = 4288085.5  # K8: 4583433.5 -> 4288085.5
$ws.Cells.Item(8, 13).Value = -4287946.5  # M8: -4583294.5 -> -4287946.5
$ws.Cells.Item(122, 8).Value = 766.1429000000001  # H122: 778.1429000000001 -> 766.1429000000001
$ws.Cells.Item(122, 9).Value = 766.1429000000001  # I122: 778.1429000000001 -> 766.1429000000001
$ws.Cells.Item(122, 11).Value = 2298.4287  # K122: 2334.4287 -> 2298.4287
$ws.Cells.Item(122, 13).Value = 151.5712999999996  # M122: 115.5712999999996 -> 151.5712999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1246.9412  # H136: 1121.225 -> 1246.9412
$ws.Cells.Item(136, 9).Value = 1045.6154  # I136: 944.4838999999999 -> 1045.6154
$ws.Cells.Item(136, 10).Value = 1901.25  # J136: 1730 -> 1901.25
$ws.Cells.Item(136, 11).Value = 3136.8462  # K136: 2833.4517 -> 3136.8462
$ws.Cells.Item(136, 12).Value = 5703.75  # L136: 5190 -> 5703.75
$ws.Cells.Item(136, 13).Value = -586.8462  # M136: -283.4516999999996 -> -586.8462
$ws.Cells.Item(136, 14).Value = -10803.75  # N136: -10290 -> -10803.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2345.3462  # H132: 2527.1667 -> 2345.3462
$ws.Cells.Item(132, 9).Value = 2374.634  # I132: 2613.6758 -> 2374.634
$ws.Cells.Item(132, 11).Value = 7123.902  # K132: 7841.0274 -> 7123.902
$ws.Cells.Item(132, 13).Value = -4593.902  # M132: -5311.0274 -> -4593.902
$ws.Cells.Item(136, 8).Value = 680.4  # H136: 660.7358400000001 -> 680.4
$ws.Cells.Item(136, 9).Value = 390.20514  # I136: 390.65854 -> 390.20514
$ws.Cells.Item(136, 10).Value = 1709.2727  # J136: 1583.5 -> 1709.2727
$ws.Cells.Item(136, 11).Value = 1170.61542  # K136: 1171.97562 -> 1170.61542
$ws.Cells.Item(136, 12).Value = 5127.8181  # L136: 4750.5 -> 5127.8181
$ws.Cells.Item(136, 13).Value = 1379.38458  # M136: 1378.02438 -> 1379.38458
$ws.Cells.Item(136, 14).Value = -10227.8181  # N136: -9850.5 -> -10227.8181
